$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 6526
$ws.Range("B1").Value = 6588
$ws.Range("C1").Value = 6638.000066928068
$ws.Range("D1").Value = 6706.531194512969
$ws.Range("E1").Value = 6526
$ws.Range("F1").Value = 6730.935737600286
$ws.Range("G1").Value = 6526
$ws.Range("H1").Value = 6526

$ws.Range("A2").Value = 6526
$ws.Range("B2").Value = 6730.935737600286
$ws.Range("C2").Value = 6526
$ws.Range("D2").Value = 6730.935737600286
$ws.Range("E2").Value = 6526
$ws.Range("F2").Value = 6526
$ws.Range("G2").Value = 6526
$ws.Range("H2").Value = 6526

$ws.Range("A3").Value = 6526
$ws.Range("B3").Value = 6730.935737600286
$ws.Range("C3").Value = 6526
$ws.Range("D3").Value = 6597.999999999999
$ws.Range("E3").Value = 6526
$ws.Range("F3").Value = 6730.935737600286
$ws.Range("G3").Value = 6526
$ws.Range("H3").Value = 6526

$ws.Range("A4").Value = 6526
$ws.Range("B4").Value = 6526
$ws.Range("C4").Value = 6526
$ws.Range("D4").Value = 6730.935737600286
$ws.Range("E4").Value = 6526
$ws.Range("F4").Value = 6526
$ws.Range("G4").Value = 6526
$ws.Range("H4").Value = 6730.935737600286

$ws.Range("A5").Value = 6526
$ws.Range("B5").Value = 6526
$ws.Range("C5").Value = 6526
$ws.Range("D5").Value = 6526
$ws.Range("E5").Value = 6607
$ws.Range("F5").Value = 6526
$ws.Range("G5").Value = 6526
$ws.Range("H5").Value = 6730.935737600286

$ws.Range("A6").Value = 6526
$ws.Range("B6").Value = 6706.531194512969
$ws.Range("C6").Value = 6526
$ws.Range("D6").Value = 6526
$ws.Range("E6").Value = 6730.935737600286
$ws.Range("F6").Value = 6526
$ws.Range("G6").Value = 6526
$ws.Range("H6").Value = 6526
